# Weekly update: insert two new "Cilantro" price records for
# Terminal Hortofrutícola Agro Chillán (week of 2022-07-29) above the
# existing data block, pushing the remaining rows (old 32-62) down to
# (34-64) and extending the used range to A1:R64.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at the top of the data block (old row 32).
$ws.Rows.Item(32).Insert()
$ws.Rows.Item(32).Insert()

# New row 32: Cilantro, Primera, week of 2022-07-29.
$ws.Range("A32").Value = 7
$ws.Range("B32").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C32").Value = "Ñuble"
$ws.Range("D32").Value = 44771
$ws.Range("E32").Value = 16
$ws.Range("F32").Value = 100112040
$ws.Range("G32").Value = "Cilantro"
$ws.Range("H32").Value = "Sin especificar"
$ws.Range("I32").Value = "Primera"
$ws.Range("J32").Value = 200
$ws.Range("K32").Value = 700
$ws.Range("L32").Value = 800
$ws.Range("M32").Value = 750
$ws.Range("N32").Value = "$/atado 0,5 a 1 kilo"
$ws.Range("O32").Value = "Provincia de Diguillín"
$ws.Range("P32").Value = 750
$ws.Range("Q32").Value = 1
$ws.Range("R32").Value = "Hortaliza"

# New row 33: Cilantro, Segunda, week of 2022-07-29.
$ws.Range("A33").Value = 7
$ws.Range("B33").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C33").Value = "Ñuble"
$ws.Range("D33").Value = 44771
$ws.Range("E33").Value = 16
$ws.Range("F33").Value = 100112040
$ws.Range("G33").Value = "Cilantro"
$ws.Range("H33").Value = "Sin especificar"
$ws.Range("I33").Value = "Segunda"
$ws.Range("J33").Value = 200
$ws.Range("K33").Value = 600
$ws.Range("L33").Value = 600
$ws.Range("M33").Value = 600
$ws.Range("N33").Value = "$/atado 0,5 a 1 kilo"
$ws.Range("O33").Value = "Provincia de Diguillín"
$ws.Range("P33").Value = 600
$ws.Range("Q33").Value = 1
$ws.Range("R33").Value = "Hortaliza"
